$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: width to match target (~22.83 chars)
$ws.Columns.Item(4).ColumnWidth = 22.83203125

# Header cell D1 "Статус" - copy header style from C1 (bold, grey fill, centered)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Статус"

# D2:D7 - bordered, centered, wrap-text style (copy border+center style from A2, then add wrap)
$ws.Range("A2").Copy()
$ws.Range("D2:D7").PasteSpecial(-4122)
$ws.Range("D2:D7").WrapText = $true

# Clear any pasted values in D2:D7 (A2 held numeric 1) so the cells stay empty
$ws.Range("D2:D7").ClearContents()

# D3 gets the "X" mark
$ws.Range("D3").Value = "X"

# Row 3 needs a taller row to fit the new content per the target layout
$ws.Rows.Item(3).RowHeight = 17

# Update the active selection to match the saved view state
$ws.Range("C11").Select() | Out-Null

$excel.CutCopyMode = $false
